# Actualización automática del inventario: agrega el producto "002SMV"
# (Chip Epson T6712) como nueva fila al final de la hoja de inventario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

$ws.Cells.Item($row, 1).Value = "002SMV"
$ws.Cells.Item($row, 2).Value = "Chip Epson"
$ws.Cells.Item($row, 3).Value = "T6712"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 10
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E29-D29)*G29"
$ws.Cells.Item($row, 9).Formula = "=D29*F29"
$ws.Cells.Item($row, 10).Value = 0
